$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.272549271583557
$ws.Range("B1").Value = 2.399449586868286
$ws.Range("D1").Value = 1.378092527389526
$ws.Range("E1").Value = 0.861849308013916
